# Appends 6 new flight-arrival rows (178-183) to the "Main Data" sheet,
# mirroring the rows already present in the workbook (columns A-L, with
# K and M left blank as in the rest of the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ A = 177; B = "Monday, Jan 16"; C = "2:52 AM";  D = "5Y587";  E = "Kaunas";    F = "(KUN)"; G = "Atlas Air ";  H = "B744"; I = "(N482MC)"; J = "5:40 AM";  L = "2 hours, 48 minutes" },
    @{ A = 178; B = "Monday, Jan 16"; C = "4:20 AM";  D = "E45070"; E = "Hurghada";  F = "(HRG)"; G = "Enter Air ";  H = "B738"; I = "(SP-ESF)"; J = "4:16 AM";  L = "0 hours, -4 minutes" },
    @{ A = 179; B = "Monday, Jan 16"; C = "6:18 AM";  D = "P81988"; E = "Cologne";   F = "(CGN)"; G = "SprintAir "; H = "AT75"; I = "(SP-SPL)"; J = "6:05 AM";  L = "0 hours, -13 minutes" },
    @{ A = 180; B = "Monday, Jan 16"; C = "9:25 AM";  D = "UNKNOWN"; E = "Ostrava";  F = "(OSR)"; G = "QA Aviation "; H = "E55P"; I = "(OK-STS)"; J = "9:18 AM"; L = "0 hours, -7 minutes" },
    @{ A = 181; B = "Monday, Jan 16"; C = "10:05 AM"; D = "LH1388"; E = "Frankfurt"; F = "(FRA)"; G = "Lufthansa "; H = "CRJ9"; I = "(D-ACNK)"; J = "10:08 AM"; L = "0 hours, 3 minutes" },
    @{ A = 182; B = "Monday, Jan 16"; C = "10:15 AM"; D = "FR3593"; E = "Milan";     F = "(BGY)"; G = "Ryanair ";   H = "B738"; I = "(SP-RSM)"; J = "10:05 AM"; L = "0 hours, -10 minutes" }
)

$startRow = 178
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 9).Value = $data.I
    $ws.Cells.Item($r, 10).Value = $data.J
    $ws.Cells.Item($r, 12).Value = $data.L
}
